# Remove duplicates, convert to xlsx and clear custom styles
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

# Row 16: split the long "Minimum load capacity/maximum capactiy(%)" text out of I16
# into a new cell P16, and replace I16 with a proper per-row description.
$ws.Range("P16").Value = "Minimum load capacity/maximum capactiy(%)"
$ws.Range("I16").Value = "Minimum stable operating level of online capacity."

# Row 19: split the long fraction-of-load-range text out of I19 into a new
# cell O19, and replace I19 with a proper per-row description.
$rq = [char]0x2019
$ws.Range("O19").Value = "fraction of the feasible load range`nabove the minimum operating level (given by`nACT_LOSPL(r,v,p,${rq}LO${rq}) ), below which the efficiency`nlosses are assumed to occur. Default value is 0.6. Unit:`nfraction of installed capacity"
$ws.Range("I19").Value = "Load level with no partial load efficiency loss "

# Reset the selection to match the post-edit state observed in the diff
$ws.Range("I19").Select()
